$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14: 3rd "-ish" measure (D14) bumped from "0H" to "0.5H"
$ws.Range("D14").Value = "0.5H"

# New row 15 for the "Greatest Common Divisor" exercise — copy row 14's
# formatting first so the new row matches the existing styles exactly.
$ws.Range("A14:E14").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A15").Value = "2/17/2013"
$ws.Range("B15").Value = "0H"
$ws.Range("C15").Value = "0.75H"
$ws.Range("D15").Value = "0.5H"
$ws.Range("E15").Value = "Greatest Common Divisor"

# Match the saved view state: scrolled back to the top, active cell on
# the row right below the new data.
$ws.Range("E16").Select() | Out-Null
